# Update cryptocurrency price/volume figures (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.626.24"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.826.26"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'309.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.4663"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.37%  "
$ws.Range("D8").Value = "'0.3597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.07132"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "'0.9025"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").Value = "'0.07707"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.829.60"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("D14").Value = "'5.269"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "'6.359"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'0.000008559"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "26.667.74"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").Value = "'14.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "'5.024"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "'10.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'1.902"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("D25").Value = "'153.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").Value = "'17.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "'1.986"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").Value = "'113.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("D29").Value = "'4.872"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'0.08814"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").Value = "'3.131"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").Value = "'2.843"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.17%  "
$ws.Range("D33").Value = "'1.169"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.86%  "
$ws.Range("D34").Value = "'0.7350"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").Value = "'4.440"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'1.080"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").Value = "'0.01930"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'0.05162"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "'2.912"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").Value = "'6.876"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "'0.5056"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'0.1497"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'8.061"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "'0.4663"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "'10.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("D47").Value = "'97.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'0.06043"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "'63.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "'35.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
